# Update attendance calculation logic results in the two report sheets.
$wb = $excel.ActiveWorkbook
$wsIndividual = $wb.Worksheets.Item("Individual Attendance")
$wsTeam = $wb.Worksheets.Item("Team Summary")

# --- Individual Attendance sheet -------------------------------------------------

# Remove the two trailing rows (Muhammad Fuaad Zameer, Shafqat Farhan) that no
# longer qualify as meetings/attendees under the new validation rule.
# Delete bottom-up so row indices of the remaining rows are not disturbed.
$wsIndividual.Rows.Item(7).EntireRow.Delete()
$wsIndividual.Rows.Item(6).EntireRow.Delete()

# Recomputed attendance figures (meetings now only count those with >= 2 real participants).
$wsIndividual.Range("C2").Value = 18
$wsIndividual.Range("D2").Value = 18
$wsIndividual.Range("E2").Value = 100

$wsIndividual.Range("C3").Value = 17
$wsIndividual.Range("D3").Value = 18
$wsIndividual.Range("E3").Value = 94.44

$wsIndividual.Range("C4").Value = 16
$wsIndividual.Range("D4").Value = 18
$wsIndividual.Range("E4").Value = 88.89

$wsIndividual.Range("C5").Value = 14
$wsIndividual.Range("D5").Value = 18
$wsIndividual.Range("E5").Value = 77.78

# Narrow the Email column slightly to match the refreshed export.
$wsIndividual.Columns.Item(2).ColumnWidth = 26.1666666666667

# --- Team Summary sheet -----------------------------------------------------------

$wsTeam.Range("B2").Value = 18
$wsTeam.Range("B3").Value = 4

# B4 must stay a literal text string ("90.28%"), not an auto-converted percentage
# number. Forcing a Text number format first stops Excel's autoconvert, then
# clearing the format afterwards drops the now-unneeded style so the cell keeps
# the workbook's default (unstyled) formatting, matching the original layout.
$wsTeam.Range("B4").NumberFormat = "@"
$wsTeam.Range("B4").Value = "90.28%"
$wsTeam.Range("B4").ClearFormats()

$wsTeam.Range("B5").Value = 3.61
